$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mobiles")

$ws.Range("G1").Value = "Status"
$ws.Range("G3").Value = "Yes"
$ws.Range("G4").Value = "Yes"
$ws.Range("G2").Value = "No"

$ws.Range("G3").Select()
